$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark review points as Accepted/Rejected (Acceptance column = E)
$ws.Range("E2").Value = "Accepted"
$ws.Range("E3").Value = "Accepted"
$ws.Range("E4").Value = "Accepted"
$ws.Range("E5").Value = "Rejected"
$ws.Range("E6").Value = "Accepted"
$ws.Range("E7").Value = "Accepted"

# Add a comment explaining the rejection of point in row 5 (Comment column = G)
$ws.Range("G5").Value = "Do Math Control is responsible of generating the string that contains the operands, operation and result. Display Control can only accept string signals"
